# costs.xlsx edit: add new expense rows to the "Έξοδα" table (Πίνακας1).
#  1) a new row is inserted right before the current row 150 (shifting every
#     following row down by one) with a new "driving license" expense;
#  2) five more new rows are appended at the end of the existing table.
# The workbook's table auto-totals (SUBTOTAL) row follows along.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$lo = $ws.ListObjects.Item(1)

# ---------------------------------------------------------------------
# 1) Insert a brand-new data row before sheet row 150.
# ---------------------------------------------------------------------
$ws.Rows.Item(150).Insert()

# Re-attach the table definition to the now one-row-taller range so the
# totals row / autofilter follow the shift (Insert() alone does not retarget
# the ListObject in this host).
$lastRow = $lo.Range.Rows.Count + 1
$lo.Resize($ws.Range($ws.Cells.Item(1,1), $ws.Cells.Item($lastRow, 10)))

$ws.Range("A150").Value = 150
$ws.Range("B150").Value = "υποχρεώσεις"
$ws.Range("C150").Value = 43211
$ws.Range("D150").Value = "δίπλωμα οδήγησης"
$ws.Range("E150").Value = "ερωτευμένος"
$ws.Range("F150").Value = 5
$ws.Range("H150").Value = "Περιστέρι"

# ---------------------------------------------------------------------
# 2) Append five more rows of new data at the bottom of the table
#    (rows 241-245, just above the existing totals row).
# ---------------------------------------------------------------------

# Row 241 - breakfast
$ws.Range("A156").Copy()
$ws.Range("A241").PasteSpecial(-4122)
$ws.Range("C151").Copy()
$ws.Range("C241").PasteSpecial(-4122)
$ws.Range("A241").Value = 2.1
$ws.Range("B241").Value = "φαγητό"
$ws.Range("C241").Value = 43263
$ws.Range("D241").Value = "πρωινό moccachino κουλούρι"
$ws.Range("E241").Value = "ελεύθερος"
$ws.Range("F241").Value = 5
$ws.Range("H241").Value = "Αμπελόκηποι"

# Row 242 - driving license
$ws.Range("A155").Copy()
$ws.Range("A242").PasteSpecial(-4122)
$ws.Range("C151").Copy()
$ws.Range("C242").PasteSpecial(-4122)
$ws.Range("A242").Value = 150
$ws.Range("B242").Value = "υποχρεώσεις"
$ws.Range("C242").Value = 43263
$ws.Range("D242").Value = "δίπλωμα οδήγησης"
$ws.Range("E242").Value = "ελεύθερος"
$ws.Range("F242").Value = 5
$ws.Range("H242").Value = "Περιστέρι"

# Row 243 - writing tablet
$ws.Range("A156").Copy()
$ws.Range("A243").PasteSpecial(-4122)
$ws.Range("C151").Copy()
$ws.Range("C243").PasteSpecial(-4122)
$ws.Range("A243").Value = 12.7
$ws.Range("B243").Value = "υλικά αγαθά"
$ws.Range("C243").Value = 43264
$ws.Range("D243").Value = "writing tablet"
$ws.Range("E243").Value = "ελεύθερος"
$ws.Range("F243").Value = 4
$ws.Range("H243").Value = "Χαλάνδρι"

# Row 244 - shoes
$ws.Range("A156").Copy()
$ws.Range("A244").PasteSpecial(-4122)
$ws.Range("C151").Copy()
$ws.Range("C244").PasteSpecial(-4122)
$ws.Range("A244").Value = 89.9
$ws.Range("B244").Value = "ένδυση"
$ws.Range("C244").Value = 43264
$ws.Range("D244").Value = "παπούτσια addidas"
$ws.Range("E244").Value = "ελεύθερος"
$ws.Range("F244").Value = 4
$ws.Range("H244").Value = "Περιστέρι"

# Row 245 - souvlaki
$ws.Range("A156").Copy()
$ws.Range("A245").PasteSpecial(-4122)
$ws.Range("C151").Copy()
$ws.Range("C245").PasteSpecial(-4122)
$ws.Range("A245").Value = 4.4
$ws.Range("B245").Value = "φαγητό"
$ws.Range("C245").Value = 43264
$ws.Range("D245").Value = "σουβλάκια σκέτα (διατροφή)"
$ws.Range("E245").Value = "ελεύθερος"
$ws.Range("F245").Value = 4
$ws.Range("H245").Value = "Περιστέρι"

$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3) Restore the selection / active cell like the author left it.
# ---------------------------------------------------------------------
$ws.Range("H151").Select()
